# Workbook: /tmp/work/before.xlsx
# Sheets: "NewLoanInput" (1), "Summary" (2), "Repayment Schedule" (3)
#
# Changes to apply (per the target diff):
#  1. Active tab moves from "Repayment Schedule" (index 2, 0-based) to
#     "NewLoanInput" (index 0, 0-based) -> activate NewLoanInput and select B20.
#  2. "Repayment Schedule" sheet loses tabSelected and its cached selection
#     moves from L23 to H20.
#  3. On "Repayment Schedule", a bunch of placeholder "0" values are cleared
#     out (Paid Date / blank-heading / various other columns for row 2, and
#     Paid Date + blank-heading column for rows 3-14).

$wb = $excel.ActiveWorkbook

$wsLoan = $wb.Worksheets.Item("NewLoanInput")
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# --- Clear the stray "0" placeholder values on the Repayment Schedule sheet ---

# Row 2: #, Days, Paid Date, (blank heading), Principal Due, Interest,
#        Penalties, In Advance, Late, (blank heading), Outstanding
$row2Cells = @("A2", "B2", "D2", "E2", "F2", "H2", "J2", "M2", "N2", "O2", "P2")
foreach ($cellRef in $row2Cells) {
    $wsSchedule.Range($cellRef).ClearContents()
}

# Rows 3-14: Paid Date (D) and the two blank-heading columns (E, O)
$wsSchedule.Range("D3:E14").ClearContents()
$wsSchedule.Range("O3:O14").ClearContents()

# --- Update selections / active sheet ---

# Repayment Schedule keeps its own cached selection (now H20) but is no
# longer the tab that is showing.
[void]$wsSchedule.Range("H20").Select()

# NewLoanInput becomes the active/visible tab with a new selection.
[void]$wsLoan.Activate()
[void]$wsLoan.Range("B20").Select()
